$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2910
$ws.Range("I62").Value = 2899.4443
$ws.Range("J62").Value = 3005
$ws.Range("K62").Value = 2899.4443
$ws.Range("L62").Value = 3005
$ws.Range("M62").Value = -2275.4443
$ws.Range("N62").Value = -4253
$ws.Range("H65").Value = 2910
$ws.Range("I65").Value = 2899.4443
$ws.Range("J65").Value = 3005
$ws.Range("K65").Value = 14497.2215
$ws.Range("L65").Value = 15025
$ws.Range("M65").Value = -11377.2215
$ws.Range("N65").Value = -21265
$ws.Range("H138").Value = 3795.2163
$ws.Range("I138").Value = 1022.5278
$ws.Range("J138").Value = 6421.9736
$ws.Range("K138").Value = 3067.5834
$ws.Range("L138").Value = 19265.9208
$ws.Range("M138").Value = 2072.4166
$ws.Range("N138").Value = -29545.9208

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19032.844
$ws.Range("I32").Value = 21192.46
$ws.Range("K32").Value = 21192.46
$ws.Range("M32").Value = -20905.46
$ws.Range("H45").Value = 1558.5172
$ws.Range("I45").Value = 1521.7
$ws.Range("J45").Value = 1640.3334
$ws.Range("K45").Value = 1521.7
$ws.Range("L45").Value = 1640.3334
$ws.Range("M45").Value = -1144.7
$ws.Range("N45").Value = -2394.3334
$ws.Range("H74").Value = 4519.6484
$ws.Range("I74").Value = 1965.1562
$ws.Range("J74").Value = 20868.4
$ws.Range("K74").Value = 1965.1562
$ws.Range("L74").Value = 20868.4
$ws.Range("M74").Value = -1091.1562
$ws.Range("N74").Value = -22616.4
$ws.Range("H77").Value = 4519.6484
$ws.Range("I77").Value = 1965.1562
$ws.Range("J77").Value = 20868.4
$ws.Range("K77").Value = 9825.780999999999
$ws.Range("L77").Value = 104342
$ws.Range("M77").Value = -5457.780999999999
$ws.Range("N77").Value = -113078
$ws.Range("H135").Value = 35107.125
$ws.Range("J135").Value = 35107.125
$ws.Range("L135").Value = 35107.125
$ws.Range("N135").Value = -45247.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1807.0444
$ws.Range("I86").Value = 1620.8636
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 1620.8636
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -497.8635999999999
$ws.Range("N86").Value = -12245
$ws.Range("H89").Value = 1807.0444
$ws.Range("I89").Value = 1620.8636
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 8104.317999999999
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -2488.317999999999
$ws.Range("N89").Value = -61227
$ws.Range("H99").Value = 975
$ws.Range("I99").Value = 975
$ws.Range("K99").Value = 975
$ws.Range("M99").Value = 523

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1416.6666
$ws.Range("I16").Value = 833.3333
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 833.3333
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -546.3333
$ws.Range("N16").Value = -2574
$ws.Range("H58").Value = 1467972.2
$ws.Range("I58").Value = 3031355.5
$ws.Range("J58").Value = 2300.5
$ws.Range("K58").Value = 3031355.5
$ws.Range("L58").Value = 2300.5
$ws.Range("M58").Value = -3031152.5
$ws.Range("N58").Value = -2706.5
$ws.Range("H94").Value = 1163.1111
$ws.Range("I94").Value = 924
$ws.Range("J94").Value = 1210.9333
$ws.Range("K94").Value = 924
$ws.Range("L94").Value = 1210.9333
$ws.Range("M94").Value = -473
$ws.Range("N94").Value = -2112.9333
$ws.Range("H99").Value = 1716.7368
$ws.Range("I99").Value = 1220.7142
$ws.Range("J99").Value = 3105.6
$ws.Range("K99").Value = 1220.7142
$ws.Range("L99").Value = 3105.6
$ws.Range("M99").Value = 277.2858000000001
$ws.Range("N99").Value = -6101.6
$ws.Range("H105").Value = 2898.5
$ws.Range("I105").Value = 1597.25
$ws.Range("J105").Value = 4199.75
$ws.Range("K105").Value = 1597.25
$ws.Range("L105").Value = 4199.75
$ws.Range("M105").Value = 149.75
$ws.Range("N105").Value = -7693.75
$ws.Range("H113").Value = 1416.6666
$ws.Range("I113").Value = 833.3333
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 833.3333
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1336.6667
$ws.Range("N113").Value = -6340
$ws.Range("H126").Value = 1716.7368
$ws.Range("I126").Value = 1220.7142
$ws.Range("J126").Value = 3105.6
$ws.Range("K126").Value = 3662.1426
$ws.Range("L126").Value = 9316.799999999999
$ws.Range("M126").Value = -1192.1426
$ws.Range("N126").Value = -14256.8
$ws.Range("H132").Value = 3487.6482
$ws.Range("I132").Value = 3623.875
$ws.Range("K132").Value = 10871.625
$ws.Range("M132").Value = -8341.625
$ws.Range("H136").Value = 1467972.2
$ws.Range("I136").Value = 3031355.5
$ws.Range("J136").Value = 2300.5
$ws.Range("K136").Value = 9094066.5
$ws.Range("L136").Value = 6901.5
$ws.Range("M136").Value = -9091516.5
$ws.Range("N136").Value = -12001.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 10423430
$ws.Range("I5").Value = 596.1111
$ws.Range("J5").Value = 23824216
$ws.Range("K5").Value = 1788.3333
$ws.Range("L5").Value = 71472648
$ws.Range("M5").Value = -1676.3333
$ws.Range("N5").Value = -71472872
$ws.Range("H104").Value = 1333.375
$ws.Range("J104").Value = 2171.75
$ws.Range("L104").Value = 6515.25
$ws.Range("N104").Value = -11757.25
$ws.Range("H107").Value = 948.5833
$ws.Range("J107").Value = 1347.091
$ws.Range("L107").Value = 4041.273
$ws.Range("N107").Value = -7881.272999999999
$ws.Range("H135").Value = 10423430
$ws.Range("I135").Value = 596.1111
$ws.Range("J135").Value = 23824216
$ws.Range("K135").Value = 5364.9999
$ws.Range("L135").Value = 214417944
$ws.Range("M135").Value = -2829.9999
$ws.Range("N135").Value = -214423014

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 34.6
$ws.Range("I2").Value = 33.333332
$ws.Range("J2").Value = 36.5
$ws.Range("K2").Value = 33.333332
$ws.Range("L2").Value = 36.5
$ws.Range("M2").Value = 79.666668
$ws.Range("N2").Value = -262.5
$ws.Range("H80").Value = 4366.7188
$ws.Range("I80").Value = 3818.7058
$ws.Range("J80").Value = 4987.8
$ws.Range("K80").Value = 3818.7058
$ws.Range("L80").Value = 4987.8
$ws.Range("M80").Value = -2820.7058
$ws.Range("N80").Value = -6983.8
$ws.Range("H83").Value = 4366.7188
$ws.Range("I83").Value = 3818.7058
$ws.Range("J83").Value = 4987.8
$ws.Range("K83").Value = 19093.529
$ws.Range("L83").Value = 24939
$ws.Range("M83").Value = -14101.529
$ws.Range("N83").Value = -34923
$ws.Range("H113").Value = 3999.5
$ws.Range("I113").Value = 3999.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3999.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1829.5
$ws.Range("N113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3141.6956
$ws.Range("I40").Value = 2720.2666
$ws.Range("J40").Value = 3931.875
$ws.Range("K40").Value = 2720.2666
$ws.Range("L40").Value = 3931.875
$ws.Range("M40").Value = -2584.2666
$ws.Range("N40").Value = -4203.875
$ws.Range("H46").Value = 2900.25
$ws.Range("I46").Value = 2967
$ws.Range("J46").Value = 2700
$ws.Range("K46").Value = 2967
$ws.Range("L46").Value = 2700
$ws.Range("M46").Value = -2779
$ws.Range("N46").Value = -3076

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1274.6316
$ws.Range("I126").Value = 1263.625
$ws.Range("K126").Value = 3790.875
$ws.Range("M126").Value = -1320.875
$ws.Range("H132").Value = 1647.0638
$ws.Range("I132").Value = 837.8484999999999
$ws.Range("J132").Value = 3554.5
$ws.Range("K132").Value = 2513.5455
$ws.Range("L132").Value = 10663.5
$ws.Range("M132").Value = 16.45450000000028
$ws.Range("N132").Value = -15723.5
